# Daily cryptos-list refresh (GitHub Actions): update Price (D) / Volume(1h) (E)
# columns for each row, plus the Cronos/Mantle row swap at rows 48-49.
# NumberFormat = "@" is applied first on any Price value that would otherwise
# be auto-parsed by Excel as a number (e.g. "1.003", "0.06741"), so it is
# stored as literal text exactly like the other price strings such as
# "28.935.44" which already fail numeric parsing on their own.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.935.44"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "1.813.88"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.44"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2746"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06741"
$ws.Range("E9").Value = "  -4.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.87"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "1.873.94"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.671"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6235"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009274"
$ws.Range("E15").Value = "  -7.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "74.32"
$ws.Range("E16").Value = "  -6.41%  "
$ws.Range("D17").Value = "28.714.30"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.434"
$ws.Range("E18").Value = "  -8.81%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.24"
$ws.Range("E20").Value = "  -9.07%  "
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.764"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.69"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1271"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.791"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06326"
$ws.Range("E28").Value = "  -6.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.403"
$ws.Range("E29").Value = "  -5.23%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.725"
$ws.Range("E31").Value = "  -3.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.682"
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("E34").Value = "  -7.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6326"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.732"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.421"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01689"
$ws.Range("E39").Value = "  -4.28%  "
$ws.Range("D40").Value = "1.131.59"
$ws.Range("E40").Value = "  -8.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8686"
$ws.Range("E41").Value = "  -5.71%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "1.971.70"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.01"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.35"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4513"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05461"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.243"
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("E51").Value = "  -0.06%  "
